# Applies the "nmv 12 06 2024" revision:
#  1. Prepends a brand-new dated correction block (title, subtitle,
#     a 2-row table, a "========" separator, and a blank paragraph)
#     before the existing first section in the document body.
#  2. Marks the (now) third table's header ("Section, Paragraph") run
#     with a <w:lastRenderedPageBreak/> as Word does when a page
#     boundary falls right before that run after the new content
#     pushes the layout down.

$d = $word.ActiveDocument

# --- Step 1: insert the new "30th June 2024" block at the very top ---
$newBlockXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">TS Pada Paatam – TS 2.4 Tamil corrections – Observed till </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr><w:t>30th June 2024</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>(ignore those which are already incorporated in your book’s version and date)</w:t></w:r></w:p><w:tbl><w:tblPr><w:tblW w:w="13301" w:type="dxa"/><w:tblInd w:w="-318" w:type="dxa"/><w:tblBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="000000"/><w:left w:val="single" w:sz="4" w:space="0" w:color="000000"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="000000"/><w:right w:val="single" w:sz="4" w:space="0" w:color="000000"/><w:insideH w:val="single" w:sz="4" w:space="0" w:color="000000"/><w:insideV w:val="single" w:sz="4" w:space="0" w:color="000000"/></w:tblBorders><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="2978"/><w:gridCol w:w="5103"/><w:gridCol w:w="5220"/></w:tblGrid><w:tr><w:trPr><w:trHeight w:val="855"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="2978" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="000000"/><w:left w:val="single" w:sz="4" w:space="0" w:color="000000"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="000000"/><w:right w:val="single" w:sz="4" w:space="0" w:color="000000"/></w:tcBorders><w:hideMark/></w:tcPr><w:p><w:pPr><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:bidi="ta-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:bidi="ta-IN"/></w:rPr><w:t>Section, Paragraph</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:bidi="ta-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:bidi="ta-IN"/></w:rPr><w:t>Reference</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5103" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="000000"/><w:left w:val="single" w:sz="4" w:space="0" w:color="000000"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="000000"/><w:right w:val="single" w:sz="4" w:space="0" w:color="000000"/></w:tcBorders><w:hideMark/></w:tcPr><w:p><w:pPr><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:bidi="ta-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:bidi="ta-IN"/></w:rPr><w:t>As Printed</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5220" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="000000"/><w:left w:val="single" w:sz="4" w:space="0" w:color="000000"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="000000"/><w:right w:val="single" w:sz="4" w:space="0" w:color="000000"/></w:tcBorders><w:hideMark/></w:tcPr><w:p><w:pPr><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:ind w:right="-18"/><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:bidi="ta-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:bidi="ta-IN"/></w:rPr><w:t>To be read as or corrected as</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="864"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="2978" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="000000"/><w:left w:val="single" w:sz="4" w:space="0" w:color="000000"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="000000"/><w:right w:val="single" w:sz="4" w:space="0" w:color="000000"/></w:tcBorders><w:hideMark/></w:tcPr><w:p><w:pPr><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:cs="Latha"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="ta-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="Latha"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="ta-IN"/></w:rPr><w:t>None</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5103" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="000000"/><w:left w:val="single" w:sz="4" w:space="0" w:color="000000"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="000000"/><w:right w:val="single" w:sz="4" w:space="0" w:color="000000"/></w:tcBorders><w:hideMark/></w:tcPr><w:p><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="BRH Malayalam Extra" w:hAnsi="BRH Malayalam Extra" w:cs="BRH Malayalam Extra"/><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:bidi="ta-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="Latha"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="ta-IN"/></w:rPr><w:t>None</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5220" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="000000"/><w:left w:val="single" w:sz="4" w:space="0" w:color="000000"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="000000"/><w:right w:val="single" w:sz="4" w:space="0" w:color="000000"/></w:tcBorders><w:hideMark/></w:tcPr><w:p><w:pPr><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="BRH Malayalam Extra" w:hAnsi="BRH Malayalam Extra" w:cs="BRH Malayalam Extra"/><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:bidi="ta-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="Latha"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="ta-IN"/></w:rPr><w:t>None</w:t></w:r></w:p></w:tc></w:tr></w:tbl><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr><w:t>========</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$startRange = $d.Range(0, 0)
$startRange.InsertXML($newBlockXml)

# --- Step 2: stamp lastRenderedPageBreak on the 3rd table's header cell ---
$lrpbXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="38E91511" w14:textId="77777777" w:rsidR="00A1406D" w:rsidRPr="005D2B3F" w:rsidRDefault="00A1406D" w:rsidP="00706EDB"><w:pPr><w:spacing w:before="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r w:rsidRPr="005D2B3F"><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:lastRenderedPageBreak/><w:t>Section, Paragraph</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$targetTable = $d.Tables.Item(3)
$headerPara = $targetTable.Cell(1, 1).Range.Paragraphs.Item(1)
$headerPara.Range.InsertXML($lrpbXml)

Write-Host "Applied TS 2.4 Tamil Pada Paatam corrections update (nmv 12 06 2024)."
